$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorganize "Regiao/Ano/Valor" rows: group by region (Brasil, Regiao Nordeste, Sergipe)
# covering 2007-2024, rather than the old year-major ordering covering 2007-2023.
$data = @(
    @('Brasil', '01/01/2007', 4.56),
    @('Brasil', '01/01/2008', 3.35),
    @('Brasil', '01/01/2009', 3.48),
    @('Brasil', '01/01/2010', 3.62),
    @('Brasil', '01/01/2011', 3.72),
    @('Brasil', '01/01/2012', 3.82),
    @('Brasil', '01/01/2013', 3.93),
    @('Brasil', '01/01/2014', 3.97),
    @('Brasil', '01/01/2015', 4.19),
    @('Brasil', '01/01/2016', 4.37),
    @('Brasil', '01/01/2017', 4.28),
    @('Brasil', '01/01/2018', 4.26),
    @('Brasil', '01/01/2019', 4.31),
    @('Brasil', '01/01/2020', 5.76),
    @('Brasil', '01/01/2021', 6.62),
    @('Brasil', '01/01/2022', 4.85),
    @('Brasil', '01/01/2023', 4.44),
    @('Brasil', '01/01/2024', 4.18),
    @('Região Nordeste', '01/01/2007', 3.69),
    @('Região Nordeste', '01/01/2008', 2.41),
    @('Região Nordeste', '01/01/2009', 2.61),
    @('Região Nordeste', '01/01/2010', 2.65),
    @('Região Nordeste', '01/01/2011', 2.88),
    @('Região Nordeste', '01/01/2012', 3.08),
    @('Região Nordeste', '01/01/2013', 3.24),
    @('Região Nordeste', '01/01/2014', 3.32),
    @('Região Nordeste', '01/01/2015', 3.61),
    @('Região Nordeste', '01/01/2016', 3.83),
    @('Região Nordeste', '01/01/2017', 3.75),
    @('Região Nordeste', '01/01/2018', 3.64),
    @('Região Nordeste', '01/01/2019', 3.76),
    @('Região Nordeste', '01/01/2020', 5.11),
    @('Região Nordeste', '01/01/2021', 5.46),
    @('Região Nordeste', '01/01/2022', 4.26),
    @('Região Nordeste', '01/01/2023', 3.94),
    @('Região Nordeste', '01/01/2024', 3.8),
    @('Sergipe', '01/01/2007', 3.8),
    @('Sergipe', '01/01/2008', 2.23),
    @('Sergipe', '01/01/2009', 2.9),
    @('Sergipe', '01/01/2010', 3.17),
    @('Sergipe', '01/01/2011', 3.53),
    @('Sergipe', '01/01/2012', 3.35),
    @('Sergipe', '01/01/2013', 3.68),
    @('Sergipe', '01/01/2014', 3.47),
    @('Sergipe', '01/01/2015', 4.54),
    @('Sergipe', '01/01/2016', 5.02),
    @('Sergipe', '01/01/2017', 4.93),
    @('Sergipe', '01/01/2018', 5.0),
    @('Sergipe', '01/01/2019', 5.16),
    @('Sergipe', '01/01/2020', 6.48),
    @('Sergipe', '01/01/2021', 6.71),
    @('Sergipe', '01/01/2022', 5.26),
    @('Sergipe', '01/01/2023', 4.72),
    @('Sergipe', '01/01/2024', 3.54)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    # Leading apostrophe forces text (avoids Excel re-parsing dd/mm/yyyy as a date),
    # then reset to the default style so no number-format override is persisted.
    $ws.Cells.Item($row, 2).Value = "'" + $item[1]
    $ws.Cells.Item($row, 2).Style = 'Normal'
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}
